# Update the "სოციალური პაკეტის მიმღებები" row (row 4) with corrected
# figures for years 2015-2021 (columns E:K). Files updated and bug fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ამბროლაური")

$ws.Range("E4").Value = 542
$ws.Range("F4").Value = 546
$ws.Range("G4").Value = 530
$ws.Range("H4").Value = 531
$ws.Range("I4").Value = 552
$ws.Range("J4").Value = 542
$ws.Range("K4").Value = 540
